$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Create a new style with yellow fill to highlight the updated rows
# (apply per-cell, skipping column C, which has no data in these rows)
$ws.Range("A5").Interior.Color = 65535
$ws.Range("B5").Interior.Color = 65535
$ws.Range("D5").Interior.Color = 65535
$ws.Range("A11").Interior.Color = 65535
$ws.Range("B11").Interior.Color = 65535
$ws.Range("D11").Interior.Color = 65535

# Update Tigran's UPI value (row 5, column B)
$ws.Range("B5").Value = 336592

# Update collaborator UPI for row 11 (Jemi and Shandao leaving) -> single UPI now
$ws.Range("B11").Value = 511294

# Update selection to reflect last active cell in the edit session
$ws.Range("H25").Select()
